$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1517.3969
$ws.Cells.Item(15, 9).Value = 1517.3969
$ws.Cells.Item(15, 11).Value = 4552.1907
$ws.Cells.Item(15, 13).Value = -4383.1907
$ws.Cells.Item(33, 8).Value = 613
$ws.Cells.Item(33, 9).Value = 546
$ws.Cells.Item(33, 10).Value = 747
$ws.Cells.Item(33, 11).Value = 546
$ws.Cells.Item(33, 12).Value = 747
$ws.Cells.Item(33, 13).Value = -317
$ws.Cells.Item(33, 14).Value = -1205
$ws.Cells.Item(82, 8).Value = 14810.5
$ws.Cells.Item(82, 9).Value = 12329.5
$ws.Cells.Item(82, 11).Value = 36988.5
$ws.Cells.Item(82, 13).Value = -36582.5
$ws.Cells.Item(85, 8).Value = 14810.5
$ws.Cells.Item(85, 9).Value = 12329.5
$ws.Cells.Item(85, 11).Value = 36988.5
$ws.Cells.Item(85, 13).Value = -35584.5
$ws.Cells.Item(100, 8).Value = 6048.875
$ws.Cells.Item(100, 9).Value = 5420.625
$ws.Cells.Item(100, 11).Value = 5420.625
$ws.Cells.Item(100, 13).Value = -4879.625
$ws.Cells.Item(101, 8).Value = 904
$ws.Cells.Item(101, 10).Value = 1912.5
$ws.Cells.Item(101, 12).Value = 5737.5
$ws.Cells.Item(101, 14).Value = -8981.5
$ws.Cells.Item(113, 8).Value = 2976.2778
$ws.Cells.Item(113, 9).Value = 2865.7
$ws.Cells.Item(113, 11).Value = 2865.7
$ws.Cells.Item(113, 13).Value = 388.3000000000002
$ws.Cells.Item(138, 8).Value = 4405.314
$ws.Cells.Item(138, 10).Value = 5216.5312
$ws.Cells.Item(138, 12).Value = 15649.5936
$ws.Cells.Item(138, 14).Value = -25929.5936

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 697.5217
$ws.Cells.Item(2, 9).Value = 616.9143
$ws.Cells.Item(2, 10).Value = 954
$ws.Cells.Item(2, 11).Value = 616.9143
$ws.Cells.Item(2, 12).Value = 954
$ws.Cells.Item(2, 13).Value = -503.9143
$ws.Cells.Item(2, 14).Value = -1180
$ws.Cells.Item(32, 8).Value = 8922.166999999999
$ws.Cells.Item(32, 10).Value = 7798.3335
$ws.Cells.Item(32, 12).Value = 7798.3335
$ws.Cells.Item(32, 14).Value = -8372.333500000001
$ws.Cells.Item(63, 8).Value = 1813.8
$ws.Cells.Item(63, 9).Value = 1941.579
$ws.Cells.Item(63, 10).Value = 1409.1666
$ws.Cells.Item(63, 11).Value = 1941.579
$ws.Cells.Item(63, 12).Value = 1409.1666
$ws.Cells.Item(63, 13).Value = -1255.579
$ws.Cells.Item(63, 14).Value = -2781.1666
$ws.Cells.Item(66, 8).Value = 1813.8
$ws.Cells.Item(66, 9).Value = 1941.579
$ws.Cells.Item(66, 10).Value = 1409.1666
$ws.Cells.Item(66, 11).Value = 9707.895
$ws.Cells.Item(66, 12).Value = 7045.833000000001
$ws.Cells.Item(66, 13).Value = -6275.895
$ws.Cells.Item(66, 14).Value = -13909.833
$ws.Cells.Item(74, 8).Value = 1398
$ws.Cells.Item(74, 9).Value = 1382.5385
$ws.Cells.Item(74, 11).Value = 1382.5385
$ws.Cells.Item(74, 13).Value = -508.5385000000001
$ws.Cells.Item(77, 8).Value = 1398
$ws.Cells.Item(77, 9).Value = 1382.5385
$ws.Cells.Item(77, 11).Value = 6912.692500000001
$ws.Cells.Item(77, 13).Value = -2544.692500000001
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()
$ws.Cells.Item(116, 8).Value = 697.5217
$ws.Cells.Item(116, 9).Value = 616.9143
$ws.Cells.Item(116, 10).Value = 954
$ws.Cells.Item(116, 11).Value = 616.9143
$ws.Cells.Item(116, 12).Value = 954
$ws.Cells.Item(116, 13).Value = 1677.0857
$ws.Cells.Item(116, 14).Value = -5542
$ws.Cells.Item(132, 8).Value = 5003432.5
$ws.Cells.Item(132, 9).Value = 3535.5
$ws.Cells.Item(132, 11).Value = 10606.5
$ws.Cells.Item(132, 13).Value = -8076.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 697.5217
$ws.Cells.Item(3, 9).Value = 616.9143
$ws.Cells.Item(3, 10).Value = 954
$ws.Cells.Item(3, 11).Value = 616.9143
$ws.Cells.Item(3, 12).Value = 954
$ws.Cells.Item(3, 13).Value = -502.9143
$ws.Cells.Item(3, 14).Value = -1182
$ws.Cells.Item(22, 8).Value = 1747.5
$ws.Cells.Item(22, 9).Value = 1719.4445
$ws.Cells.Item(22, 11).Value = 1719.4445
$ws.Cells.Item(22, 13).Value = -1546.4445
$ws.Cells.Item(107, 8).Value = 2632.9048
$ws.Cells.Item(107, 9).Value = 2899.8
$ws.Cells.Item(107, 11).Value = 2899.8
$ws.Cells.Item(107, 13).Value = -979.8000000000002
$ws.Cells.Item(134, 8).Value = 8334847
$ws.Cells.Item(134, 9).Value = 929.2222
$ws.Cells.Item(134, 10).Value = 33336600
$ws.Cells.Item(134, 11).Value = 2787.6666
$ws.Cells.Item(134, 12).Value = 100009800
$ws.Cells.Item(134, 13).Value = -252.6666
$ws.Cells.Item(134, 14).Value = -100014870

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(52, 8).Value = 498.25
$ws.Cells.Item(52, 10).Value = 498.25
$ws.Cells.Item(52, 12).Value = 1494.75
$ws.Cells.Item(52, 14).Value = -2026.75
$ws.Cells.Item(56, 8).Value = 16015.346
$ws.Cells.Item(56, 9).Value = 16015.346
$ws.Cells.Item(56, 11).Value = 16015.346
$ws.Cells.Item(56, 13).Value = -15485.346
$ws.Cells.Item(64, 8).Value = 7301.4546
$ws.Cells.Item(64, 9).Value = 3167.5
$ws.Cells.Item(64, 11).Value = 9502.5
$ws.Cells.Item(64, 13).Value = -9232.5
$ws.Cells.Item(67, 8).Value = 7301.4546
$ws.Cells.Item(67, 9).Value = 3167.5
$ws.Cells.Item(67, 11).Value = 9502.5
$ws.Cells.Item(67, 13).Value = -8566.5
$ws.Cells.Item(107, 8).Value = 3505209
$ws.Cells.Item(107, 9).Value = 4059.2
$ws.Cells.Item(107, 10).Value = 4338816
$ws.Cells.Item(107, 11).Value = 12177.6
$ws.Cells.Item(107, 12).Value = 13016448
$ws.Cells.Item(107, 13).Value = -10257.6
$ws.Cells.Item(107, 14).Value = -13020288
$ws.Cells.Item(138, 8).Value = 12317.421
$ws.Cells.Item(138, 9).Value = 9846.941000000001
$ws.Cells.Item(138, 11).Value = 29540.823
$ws.Cells.Item(138, 13).Value = -24400.823
$ws.Cells.Item(139, 8).Value = 4703.1284
$ws.Cells.Item(139, 9).Value = 2880.3103
$ws.Cells.Item(139, 11).Value = 8640.930899999999
$ws.Cells.Item(139, 13).Value = -3500.930899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 10599.667
$ws.Cells.Item(10, 10).Value = 10599.667
$ws.Cells.Item(10, 12).Value = 10599.667
$ws.Cells.Item(10, 14).Value = -10937.667
$ws.Cells.Item(21, 8).Value = 28594284
$ws.Cells.Item(21, 9).Value = 33355332
$ws.Cells.Item(21, 11).Value = 33355332
$ws.Cells.Item(21, 13).Value = -33355159
$ws.Cells.Item(30, 8).Value = 28594284
$ws.Cells.Item(30, 9).Value = 33355332
$ws.Cells.Item(30, 11).Value = 33355332
$ws.Cells.Item(30, 13).Value = -33355227
$ws.Cells.Item(33, 8).Value = 20000
$ws.Cells.Item(33, 10).Value = 20000
$ws.Cells.Item(33, 12).Value = 20000
$ws.Cells.Item(33, 14).Value = -20504
$ws.Cells.Item(35, 8).Value = 29153
$ws.Cells.Item(35, 9).Value = 21507.5
$ws.Cells.Item(35, 11).Value = 21507.5
$ws.Cells.Item(35, 13).Value = -21209.5
$ws.Cells.Item(97, 8).Value = 7199.35
$ws.Cells.Item(97, 10).Value = 11860.5
$ws.Cells.Item(97, 12).Value = 11860.5
$ws.Cells.Item(97, 14).Value = -12852.5
$ws.Cells.Item(102, 8).Value = 2485.7222
$ws.Cells.Item(102, 10).Value = 4249.75
$ws.Cells.Item(102, 12).Value = 4249.75
$ws.Cells.Item(102, 14).Value = -7493.75
$ws.Cells.Item(113, 8).Value = 1685766.5
$ws.Cells.Item(113, 10).Value = 2648077.2
$ws.Cells.Item(113, 12).Value = 2648077.2
$ws.Cells.Item(113, 14).Value = -2652417.2
$ws.Cells.Item(122, 8).Value = 2764.3513
$ws.Cells.Item(122, 9).Value = 1849.4062
$ws.Cells.Item(122, 11).Value = 5548.2186
$ws.Cells.Item(122, 13).Value = -3098.2186
$ws.Cells.Item(132, 8).Value = 4341694.5
$ws.Cells.Item(132, 9).Value = 2833.3684
$ws.Cells.Item(132, 11).Value = 8500.1052
$ws.Cells.Item(132, 13).Value = -5970.1052

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 4443.5
$ws.Cells.Item(25, 9).Value = 4443.5
$ws.Cells.Item(25, 11).Value = 4443.5
$ws.Cells.Item(25, 13).Value = -4213.5
$ws.Cells.Item(40, 8).Value = 3884.1875
$ws.Cells.Item(40, 9).Value = 2864.634
$ws.Cells.Item(40, 11).Value = 2864.634
$ws.Cells.Item(40, 13).Value = -2728.634
$ws.Cells.Item(122, 8).Value = 3462.8572
$ws.Cells.Item(122, 9).Value = 3306.7778
$ws.Cells.Item(122, 11).Value = 9920.3334
$ws.Cells.Item(122, 13).Value = -7470.3334
$ws.Cells.Item(131, 8).Value = 163440.33
$ws.Cells.Item(131, 10).Value = 163440.33
$ws.Cells.Item(131, 12).Value = 163440.33
$ws.Cells.Item(131, 14).Value = -173520.33
$ws.Cells.Item(138, 8).Value = 74500
$ws.Cells.Item(138, 10).Value = 99000
$ws.Cells.Item(138, 12).Value = 99000
$ws.Cells.Item(138, 14).Value = -109280

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(46, 8).Value = 125960.25
$ws.Cells.Item(46, 10).Value = 125960.25
$ws.Cells.Item(46, 12).Value = 125960.25
$ws.Cells.Item(46, 14).Value = -126422.25
$ws.Cells.Item(126, 8).Value = 3308.524
$ws.Cells.Item(126, 9).Value = 5120
$ws.Cells.Item(126, 10).Value = 893.2222
$ws.Cells.Item(126, 11).Value = 15360
$ws.Cells.Item(126, 12).Value = 2679.6666
$ws.Cells.Item(126, 13).Value = -12890
$ws.Cells.Item(126, 14).Value = -7619.6666
$ws.Cells.Item(134, 8).Value = 125960.25
$ws.Cells.Item(134, 10).Value = 125960.25
$ws.Cells.Item(134, 12).Value = 377880.75
$ws.Cells.Item(134, 14).Value = -382950.75
$ws.Cells.Item(135, 8).Value = 72064.45
$ws.Cells.Item(135, 10).Value = 72064.45
$ws.Cells.Item(135, 12).Value = 72064.45
$ws.Cells.Item(135, 14).Value = -82204.45
